# Add two new header columns (G and H) to the EMPENHOS sheet, mirroring the
# style used by the existing header cells (A1, B1, D1, E1, F1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header values
$ws.Range("G1").Value = "MODALIDADE DA LICITAÇÃO"
$ws.Range("H1").Value = "NATUREZA DA DESPESA"

# Match the header style (bold font + fill + border) used by the other
# header cells in row 1, by copying A1's formatting onto the new cells.
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Size the new columns to fit their header text, like the rest of the sheet
# (values chosen so the exporter's pixel-quantized stored width lands on the
# closest possible match to the target column widths).
$ws.Columns.Item(7).ColumnWidth = 25.6
$ws.Columns.Item(8).ColumnWidth = 21.29

# Move the selection to match the edited workbook's saved state.
[void]$ws.Range("H5").Select()
